# Applies the "all changes in test script" commit to the LoginData sheet:
#  - fills in the (previously empty) row 53 / row 55 "private tool" test rows
#  - fixes up the "search" value cells C45 / C57
#  - appends the new "private track" add/edit section as rows 62-65
#  - updates the window selection to match the author's final cursor position

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 45: private-test search box now pre-filled with "test tool private" ---
$ws.Cells.Item(45, 3).Value = "test tool private"

# --- Row 53: "Add private Leadership tools" sample data row ---
$ws.Cells.Item(53, 2).Value = "test tool private"
$ws.Cells.Item(53, 3).Value = "ttool"
$ws.Cells.Item(53, 4).Value = "Experion Whitelabel"
$ws.Cells.Item(53, 5).Value = " desc"

# --- Row 55: "Edit private leadership tools" sample data row ---
$ws.Cells.Item(55, 2).Value = "tool test private"
$ws.Cells.Item(55, 3).Value = "tool"
$ws.Cells.Item(55, 4).Value = "Experion Whitelabel"
$ws.Cells.Item(55, 5).Value = " content"

# --- Row 57: track search box now pre-filled with "private track" ---
$ws.Cells.Item(57, 3).Value = "private track"

# --- Row 62: "private add track" header (copy the highlighted-label style from B52:E52) ---
$ws.Cells.Item(52, 2).Resize(1, 4).Copy() | Out-Null
$ws.Cells.Item(62, 2).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(62, 1).Value = "private add track"
$ws.Cells.Item(62, 2).Value = "name"
$ws.Cells.Item(62, 3).Value = "slug"
$ws.Cells.Item(62, 4).Value = "list"
$ws.Cells.Item(62, 5).Value = "desc"

# --- Row 63: sample data for the "private add track" row ---
$ws.Cells.Item(63, 2).Value = "private track"
$ws.Cells.Item(63, 3).Value = "test"
$ws.Cells.Item(63, 4).Value = "Experion Whitelabel"
$ws.Cells.Item(63, 5).Value = " desc"

# --- Row 64: "edit private track" header (copy header style from row 60 "Edit Track", B:F, then G from B52) ---
$ws.Cells.Item(60, 2).Resize(1, 5).Copy() | Out-Null
$ws.Cells.Item(64, 2).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(52, 2).Copy() | Out-Null
$ws.Cells.Item(64, 7).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(64, 1).Value = "edit private track"
$ws.Cells.Item(64, 2).Value = "name "
$ws.Cells.Item(64, 3).Value = "slug"
$ws.Cells.Item(64, 4).Value = "list"
$ws.Cells.Item(64, 5).Value = "Desccription"
$ws.Cells.Item(64, 6).Value = "tags"
$ws.Cells.Item(64, 7).Value = "courses"

# --- Row 65: sample data for the "edit private track" row ---
$ws.Cells.Item(65, 2).Value = "track 1"
$ws.Cells.Item(65, 3).Value = "test 1"
$ws.Cells.Item(65, 4).Value = "Experion Whitelabel"
$ws.Cells.Item(65, 5).Value = " desc1"
$ws.Cells.Item(65, 6).Value = "ttc"
$ws.Cells.Item(65, 7).Value = "trial test course"

# --- Final cursor/selection position left by the author ---
$win = $excel.ActiveWindow
$win.ScrollRow = 40
$win.ScrollColumn = 1
$ws.Range("A66").Select() | Out-Null
